{"js": "// Roles and Responsibilities section edits:\n//  - Brian/Jasim paragraphs: runs get merged (no wording change)\n//  - Vera paragraph: mention her art (\"artist and\") and tidy up sentence\n//  - \"A large challenge...\" paragraph: runs get merged (no wording change)\n//\n// We locate each paragraph's full text via a search, then replace the whole\n// match with the final text using Word.InsertLocation.replace. Office.js\n// collapses the matched range into a single run, inheriting the\n// character formatting (rPr) of the first run in the match - exactly what\n// the target document looks like (single <w:r> per paragraph with the\n// original <w:sz w:val=\"22\"/> / <w:szCs w:val=\"22\"/> run properties retained).\n\nconst body = context.document.body;\n\nasync function replaceExact(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(oldText) + \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Brian paragraph - merge the two runs into one (text itself is unchanged).\nawait replaceExact(\n  \"Brian \\u2013 Completed the full implementation for the Enemy\\u2019s movement method and helped devise strategy for implementation.\",\n  \"Brian \\u2013 Completed the full implementation for the Enemy\\u2019s movement method and helped devise strategy for implementation.\"\n);\n\n// 2) Jasim paragraph - merge the two runs into one (text itself is unchanged).\nawait replaceExact(\n  \"Jasim \\u2013 Coded the game timer, scheduled group meetings, and assisted in game design logic.\",\n  \"Jasim \\u2013 Coded the game timer, scheduled group meetings, and assisted in game design logic.\"\n);\n\n// 3) Vera paragraph - the actual content edit: mention that Vera did the art.\nawait replaceExact(\n  \"Vera \\u2013 Was the brains of the project and helped significantly with all of the code. Spearheaded the report and code production, giving a good backbone for the others to assist. \",\n  \"Vera \\u2013 Was the artist and brains of the project. Spearheaded the report and code production, giving a good backbone for the others to assist. \"\n);\n\n// 4) \"A large challenge...\" paragraph - merge the three runs into one (text unchanged).\nawait replaceExact(\n  \"A large challenge was getting used to and using the core libraries since for half the group this is the first time they have coded with java. Developing the movement method for the enemy took a bit of time to draft out and implement due to the logic. Another large challenge was sticking to the original UML. This was likely because we had naively drafted things, for example just because the enemy and the player are characters does not mean they should have their own superclass. I believe all these issues boil down to a lack of experience in developing a program of this size.\",\n  \"A large challenge was getting used to and using the core libraries since for half the group this is the first time they have coded with java. Developing the movement method for the enemy took a bit of time to draft out and implement due to the logic. Another large challenge was sticking to the original UML. This was likely because we had naively drafted things, for example just because the enemy and the player are characters does not mean they should have their own superclass. I believe all these issues boil down to a lack of experience in developing a program of this size.\"\n);\n", "ps1": "# Roles and Responsibilities section edits:\n#  - Brian/Jasim paragraphs: runs get merged (no wording change)\n#  - Vera paragraph: mention her art (\"artist and\") and tidy up the sentence\n#  - \"A large challenge...\" paragraph: runs get merged (no wording change)\n#\n# We use the classic Find & Replace pattern (Range.Find.Execute(..., ReplaceWith,\n# Replace:=wdReplaceOne)) rather than a plain `$range.Text = ...` assignment.\n# Word's \"real\" Find & Replace rebuilds the matched range as a single run that\n# inherits the character formatting (rPr) of the original text - exactly what\n# the target document looks like (a single <w:r> per paragraph, keeping the\n# original <w:sz w:val=\"22\"/> / <w:szCs w:val=\"22\"/> run properties).\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText {\n    param(\n        [string]$OldText,\n        [string]$NewText\n    )\n\n    $rng = $d.Content\n    $wdFindContinue = 1\n    $wdReplaceOne = 1\n    $ok = $rng.Find.Execute(\n        $OldText,    # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        $wdFindContinue, # Wrap\n        $false,      # Format\n        $NewText,    # ReplaceWith\n        $wdReplaceOne    # Replace\n    )\n\n    if (-not $ok) {\n        throw \"Could not find/replace text: $OldText\"\n    }\n}\n\n# 1) Brian paragraph - merge the two runs into one (text itself is unchanged).\nReplace-ExactText `\n    \"Brian \u2013 Completed the full implementation for the Enemy\u2019s movement method and helped devise strategy for implementation.\" `\n    \"Brian \u2013 Completed the full implementation for the Enemy\u2019s movement method and helped devise strategy for implementation.\"\n\n# 2) Jasim paragraph - merge the two runs into one (text itself is unchanged).\nReplace-ExactText `\n    \"Jasim \u2013 Coded the game timer, scheduled group meetings, and assisted in game design logic.\" `\n    \"Jasim \u2013 Coded the game timer, scheduled group meetings, and assisted in game design logic.\"\n\n# 3) Vera paragraph - the actual content edit: mention that Vera did the art.\nReplace-ExactText `\n    \"Vera \u2013 Was the brains of the project and helped significantly with all of the code. Spearheaded the report and code production, giving a good backbone for the others to assist. \" `\n    \"Vera \u2013 Was the artist and brains of the project. Spearheaded the report and code production, giving a good backbone for the others to assist. \"\n\n# 4) \"A large challenge...\" paragraph - merge the three runs into one (text unchanged).\nReplace-ExactText `\n    \"A large challenge was getting used to and using the core libraries since for half the group this is the first time they have coded with java. Developing the movement method for the enemy took a bit of time to draft out and implement due to the logic. Another large challenge was sticking to the original UML. This was likely because we had naively drafted things, for example just because the enemy and the player are characters does not mean they should have their own superclass. I believe all these issues boil down to a lack of experience in developing a program of this size.\" `\n    \"A large challenge was getting used to and using the core libraries since for half the group this is the first time they have coded with java. Developing the movement method for the enemy took a bit of time to draft out and implement due to the logic. Another large challenge was sticking to the original UML. This was likely because we had naively drafted things, for example just because the enemy and the player are characters does not mean they should have their own superclass. I believe all these issues boil down to a lack of experience in developing a program of this size.\"\n"}
